{"js": "// Target change (from the supplied OOXML diff):\n//   customXml/item1.xml  <-> customXml/item2.xml <-> customXml/item3.xml <-> customXml/item4.xml\n//   customXml/itemProps1.xml .. itemProps4.xml are relabelled the same way\n//\n// Careful byte-for-byte analysis of the diff shows this is a *pure*\n// cyclic renumbering of the package's four custom-XML parts:\n//   - item1.xml (p:properties / documentManagement)      -> item4.xml\n//   - item2.xml (b:Sources bibliography)                 -> item1.xml\n//   - item3.xml (FormTemplates)                           -> item2.xml\n//   - item4.xml (ct:contentTypeSchema)                    -> item3.xml\n// and itemPropsN.xml follow the same rotation so that each datastore\n// item keeps pointing at its own item*.xml twin.\n//\n// Crucially, every one of those eight blobs (id / namespace / XML body)\n// is byte-identical before and after - only the *file name* inside the\n// zip package changed. The document's real, API-visible content (body\n// text, styles, headers/footers, the CustomXmlParts collection as seen\n// through the object model - id, namespaceUri, xml) is 100% unchanged;\n// this is a harmless repackaging artifact (most likely produced by\n// Word/SharePoint re-syncing its Document Information Panel parts),\n// not a content edit.\n//\n// This is also not something the Word JavaScript API can perform:\n//   - CustomXmlPart.setXml() always raises GeneralException (built-in\n//     Document Information Panel parts are read-only through this\n//     surface, mirroring real Word behaviour).\n//   - CustomXmlPart.insertElement/updateElement/deleteElement/\n//     insertAttribute/updateAttribute/deleteAttribute all raise too.\n//   - CustomXmlPartCollection.add() mints a brand new part with a new\n//     id - it cannot reproduce (or choose) an existing GUID, so\n//     delete()+add() cannot restore the original identities in their\n//     new slots either.\n// There is therefore no supported Office.js call that reproduces this\n// packaging-only rename without destroying/ fabricating part identities\n// that don't match the target at all - doing so would move the\n// document further from the target, not closer.\n//\n// Net effect: nothing reachable through context.document changes, so\n// this script intentionally performs no mutation. (Read-only sanity\n// check below, left in so `context.sync()` still exercises the API.)\n\nconst parts = context.document.customXmlParts;\nparts.load(\"items\");\nawait context.sync();\n", "ps1": "# Target change (from the supplied OOXML diff):\n#   customXml/item1.xml  <-> customXml/item2.xml <-> customXml/item3.xml <-> customXml/item4.xml\n#   customXml/itemProps1.xml .. itemProps4.xml are relabelled the same way\n#\n# Careful byte-for-byte analysis of the diff shows this is a *pure*\n# cyclic renumbering of the package's four custom-XML parts:\n#   - item1.xml (p:properties / documentManagement)      -> item4.xml\n#   - item2.xml (b:Sources bibliography)                 -> item1.xml\n#   - item3.xml (FormTemplates)                           -> item2.xml\n#   - item4.xml (ct:contentTypeSchema)                    -> item3.xml\n# and itemPropsN.xml follow the same rotation so that each datastore\n# item keeps pointing at its own item*.xml twin.\n#\n# Every one of those eight blobs (id / namespace / XML body) is\n# byte-identical before and after - only the *file name* inside the zip\n# package changed. The document's real, COM-visible content (body text,\n# styles, headers/footers, CustomXMLParts - Id, NamespaceURI, XML) is\n# 100% unchanged; this is a harmless repackaging artifact (most likely\n# produced by Word/SharePoint re-syncing its Document Information Panel\n# parts), not a content edit.\n#\n# It is also not something the Word COM object model can perform here:\n# $d.CustomXMLParts / SelectByID / Item(...) resolve syntactically but\n# do not expose working Id/NamespaceURI/XML getters or setters for these\n# Document-Information-Panel-style parts, and there is no \"rename part\"\n# verb in the object model at all (CustomXMLParts.Add mints a brand new\n# part/id; it cannot reproduce an existing GUID, and CustomXMLPart has\n# no Save-As/Move operation). So there is no supported COM call that\n# reproduces this packaging-only rename without fabricating part\n# identities that don't match the target - doing so would move the\n# document further from the target, not closer.\n#\n# Net effect: nothing reachable through $d changes, so this script\n# intentionally performs no mutation. (Read-only sanity check below,\n# left in so the interpreter still touches the object model.)\n\n$d = $word.ActiveDocument\n$parts = $d.CustomXMLParts\nWrite-Output (\"CustomXMLParts.Count=\" + $parts.Count)\n"}
